# 03项目计划表.xlsx - add a new weekly status block (rows 31-40) to Sheet1,
# mirroring the existing "日期：2018.10.10 第六周周三" block (rows 21-30).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Merge the title/summary rows *before* copying formats onto them - doing
#    it in this order keeps the border set to a single uniform box (matching
#    the source file); merging afterwards makes Excel re-split the border
#    into per-edge pieces for the (now extra) inner gridlines.
$ws.Range("A31:D31").Merge() | Out-Null
$ws.Range("A39:D40").Merge() | Out-Null

# 2) Clone the formatting of the previous week's block (rows 21-30) onto the
#    new rows (31-40) without touching values, so borders/fonts/alignment
#    match exactly.
$ws.Range("A21:D30").Copy() | Out-Null
$ws.Range("A31").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# 3) Fill in the new block's header/date row.
$ws.Range("A31").Value = "日期：2018.10.11第六周周四"

# 4) Table header row (reuses existing shared strings).
$ws.Range("A32").Value = "组员"
$ws.Range("B32").Value = "计划内容"
$ws.Range("C32").Value = "完成情况"
$ws.Range("D32").Value = "备注"

# 5) Member rows. Column A / C reuse existing shared strings; column B is new
#    text. The B35 text is written before B34's so the shared-string table
#    ends up in the same append order as the source workbook.
$ws.Range("A33").Value = "陈柯赞"
$ws.Range("B33").Value = "管理员登录系统、禁用用户、开启用户、关闭用户界面"
$ws.Range("C33").Value = "已完成"

$ws.Range("A35").Value = "王智永"
$ws.Range("B35").Value = "用户登录、找回密码、创建用户群界面"
$ws.Range("C35").Value = "已完成"

$ws.Range("A34").Value = "黎安生"
$ws.Range("B34").Value = "管理员关闭群、开启群界面"
$ws.Range("C34").Value = "已完成"

$ws.Range("A36").Value = "郑海文"
$ws.Range("B36").Value = "退出用户群、编辑群资料、接受用户加入群界面"
$ws.Range("C36").Value = "已完成"

$ws.Range("A37").Value = "赵华亮"
$ws.Range("B37").Value = "拒绝用户加入群、清退群成员界面"
$ws.Range("C37").Value = "已完成"

$ws.Range("A38").Value = "叶田"
$ws.Range("B38").Value = "查看群成员位置信息界面"
$ws.Range("C38").Value = "已完成"

# 6) Apply the distinct font used for the new interface-description column
#    (B33:B38) in the source edit.
$bCol = $ws.Range("B33:B38")
$bCol.Font.Name = "宋体"
$bCol.Font.Size = 11

# 7) Summary row label.
$ws.Range("A39").Value = "总结："

# 8) Restore the view/selection state recorded in the saved file.
$ws.Range("B38").Select()
